# Avg Vehicle Loading.xlsx edit:
#  - Add a new "US AVLo" sheet (between "rail" and "AVLo-passengers") containing
#    US vehicle-loading reference data (passengers + tons blocks).
#  - Point AVLo-freight's "ships" row (domestic shipping) at the new sheet's
#    "ships" tons figure via a formula, instead of a hard-coded literal, and
#    apply the matching number format.
#  - Update sheet selections and the active tab left by the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "US AVLo" worksheet right after "rail".
# ---------------------------------------------------------------------------
$railSheet = $wb.Worksheets.Item("rail")
$usAvlo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $railSheet)
$usAvlo.Name = "US AVLo"

# ---- Block 1: "Vehicle Loading (passengers)", row 1 header = years 2015-2050 (B:AK)
$usAvlo.Cells.Item(1, 1).Value = "Vehicle Loading (passengers)"
$col = 2
for ($yr = 2015; $yr -le 2050; $yr++) {
    $usAvlo.Cells.Item(1, $col).Value = $yr
    $col++
}

$passengerRows = @{
    2 = @{ Label = "LDVs";       Value = 1.67 }
    3 = @{ Label = "HDVs";       Value = 21.196137258578663 }
    4 = @{ Label = "aircraft";   Value = 111.39416306433705 }
    5 = @{ Label = "rail";       Value = 486.56731685074101 }
    6 = @{ Label = "ships";      Value = 1 }
    7 = @{ Label = "motorbikes"; Value = 1.2700756740871355 }
}
foreach ($r in 2..7) {
    $info = $passengerRows[$r]
    $usAvlo.Cells.Item($r, 1).Value = $info.Label
    for ($c = 2; $c -le 37; $c++) {
        $usAvlo.Cells.Item($r, $c).Value = $info.Value
    }
}

# ---- Block 2: "Vehicle Loading (tons)", row 9 header = years 2016-2050 (B:AJ)
$usAvlo.Cells.Item(9, 1).Value = "Vehicle Loading (tons)"
$col = 2
for ($yr = 2016; $yr -le 2050; $yr++) {
    $usAvlo.Cells.Item(9, $col).Value = $yr
    $col++
}

$tonRows = @{
    10 = @{ Label = "LDVs";       Value = 1 }
    11 = @{ Label = "HDVs";       Value = 16 }
    12 = @{ Label = "aircraft";   Value = 41.989116133258747 }
    13 = @{ Label = "rail";       Value = 3512.35916421195 }
    14 = @{ Label = "ships";      Value = 1974.4736422180429 }
    15 = @{ Label = "motorbikes"; Value = 0 }
}
foreach ($r in 10..15) {
    $info = $tonRows[$r]
    $usAvlo.Cells.Item($r, 1).Value = $info.Label
    for ($c = 2; $c -le 36; $c++) {
        $usAvlo.Cells.Item($r, $c).Value = $info.Value
    }
}

# Leave the new sheet with rows 9:15 selected (matches the authored state).
$usAvlo.Activate()
$usAvlo.Rows("9:15").Select()

# ---------------------------------------------------------------------------
# 2. AVLo-freight: drive the "ships" (domestic shipping) row off the new
#    "US AVLo" tons figure instead of the old hard-coded constant, and give
#    it the "0_ " number format used elsewhere for integer-like figures.
# ---------------------------------------------------------------------------
$avloFreight = $wb.Worksheets.Item("AVLo-freight")
for ($c = 2; $c -le 36; $c++) {
    $cell = $avloFreight.Cells.Item(6, $c)
    $cell.Formula = "='US AVLo'!`$B`$14"
    $cell.NumberFormat = "0_ "
}

$avloFreight.Activate()
$avloFreight.Range("B6:AJ6").Select()

# ---------------------------------------------------------------------------
# 3. AVLo-passengers: move the lingering selection.
# ---------------------------------------------------------------------------
$avloPassengers = $wb.Worksheets.Item("AVLo-passengers")
$avloPassengers.Activate()
$avloPassengers.Range("B7").Select()

# ---------------------------------------------------------------------------
# 4. Leave "About" as the active sheet/tab.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Activate()
